# "Create ROM 2 testbench" - append a new time-log row (row 43) to the
# Time Record sheet, mirroring the existing rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting of the previous row (42) down into row 43 so the
#     new cells pick up the same styles (date alignment, time number
#     format, formula alignment, task/notes alignment) without dragging
#     along row 42's values/formulas. Column G is left alone - this entry
#     has no Notes, so row 43 shouldn't get a G cell at all. ---
$ws.Range("A42:F42").Copy()
$ws.Range("A43:F43").PasteSpecial(-4122)

# --- Date (column A) is stored as text ("4.4.2020"), just like the rest
#     of the column, so force text entry before writing it. ---
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "4.4.2020"
# Re-apply row 42's formatting so the temporary "@" number format we used
# to force text entry doesn't leak into the saved style.
$ws.Range("A42").Copy()
$ws.Range("A43").PasteSpecial(-4122)

# --- From / To (columns B and C) are time-of-day fractions. ---
$ws.Range("B43").Value = 0.4861111111111111
$ws.Range("C43").Value = 0.48958333333333331

# --- Time (column D) continues the existing (=C-B) pattern that runs
#     down the whole sheet. ---
$ws.Range("D43").Formula = "=C43-B43"

# --- Unit / Task (columns E and F); no Notes (column G) for this entry. ---
$ws.Range("E43").Value = "ROM 2"
$ws.Range("F43").Value = "Testbench"

# --- Move the view/selection down to the new last row, matching where
#     the author was working when they added the entry. ---
$ws.Range("G43").Select()
